# Auto-generated Excel COM-interop script to apply Tonberry_Profits price/profit updates
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 493.44446
$ws.Range("I2").Value = 407.66666
$ws.Range("J2").Value = 665
$ws.Range("K2").Value = 407.66666
$ws.Range("L2").Value = 665
$ws.Range("N2").Value = -891
$ws.Range("H17").Value = 1625.322
$ws.Range("J17").Value = 1381.8948
$ws.Range("L17").Value = 4145.6844
$ws.Range("N17").Value = -4481.6844
$ws.Range("H40").Value = 1014.7659
$ws.Range("I40").Value = 995.4878
$ws.Range("K40").Value = 995.4878
$ws.Range("M40").Value = -820.4878
$ws.Range("H62").Value = 1597.8334
$ws.Range("J62").Value = 1700
$ws.Range("L62").Value = 1700
$ws.Range("N62").Value = -2948
$ws.Range("H65").Value = 1597.8334
$ws.Range("J65").Value = 1700
$ws.Range("L65").Value = 8500
$ws.Range("N65").Value = -14740
$ws.Range("H86").Value = 155452.62
$ws.Range("I86").Value = 177374.42
$ws.Range("J86").Value = 2000
$ws.Range("K86").Value = 177374.42
$ws.Range("L86").Value = 2000
$ws.Range("M86").Value = -176251.42
$ws.Range("N86").Value = -4246
$ws.Range("H89").Value = 155452.62
$ws.Range("I89").Value = 177374.42
$ws.Range("J89").Value = 2000
$ws.Range("K89").Value = 886872.1000000001
$ws.Range("L89").Value = 10000
$ws.Range("M89").Value = -881256.1000000001
$ws.Range("N89").Value = -21232
$ws.Range("H100").Value = 1179.3
$ws.Range("I100").Value = 925.375
$ws.Range("K100").Value = 925.375
$ws.Range("M100").Value = -384.375
$ws.Range("H107").Value = 390.42856
$ws.Range("J107").Value = 2150
$ws.Range("L107").Value = 2150
$ws.Range("N107").Value = -5990
$ws.Range("H116").Value = 16769.9
$ws.Range("I116").Value = 26924.75
$ws.Range("J116").Value = 10000
$ws.Range("K116").Value = 26924.75
$ws.Range("L116").Value = 10000
$ws.Range("M116").Value = -23482.75
$ws.Range("N116").Value = -16884
$ws.Range("H135").Value = 35714724
$ws.Range("I135").Value = 474.32
$ws.Range("K135").Value = 4268.88
$ws.Range("M135").Value = -1733.88
$ws.Range("H137").Value = 957.5
$ws.Range("I137").Value = 627.04
$ws.Range("J137").Value = 1875.4445
$ws.Range("K137").Value = 1881.12
$ws.Range("L137").Value = 5626.333500000001
$ws.Range("M137").Value = 668.8800000000001
$ws.Range("N137").Value = -10726.3335
$ws.Range("H138").Value = 1565.039
$ws.Range("I138").Value = 1231.9344
$ws.Range("J138").Value = 2835
$ws.Range("K138").Value = 3695.8032
$ws.Range("L138").Value = 8505
$ws.Range("M138").Value = 1444.1968
$ws.Range("N138").Value = -18785
$ws.Range("H141").Value = 850005.6
$ws.Range("I141").Value = 1000806.7
$ws.Range("K141").Value = 3002420.1
$ws.Range("M141").Value = -2997240.1

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1389749.8
$ws.Range("I2").Value = 1852333
$ws.Range("K2").Value = 1852333
$ws.Range("M2").Value = -1852220
$ws.Range("H32").Value = 3173.2825
$ws.Range("I32").Value = 2519.2727
$ws.Range("J32").Value = 9468.125
$ws.Range("K32").Value = 2519.2727
$ws.Range("L32").Value = 9468.125
$ws.Range("M32").Value = -2232.2727
$ws.Range("N32").Value = -10042.125
$ws.Range("H45").Value = 1689.6522
$ws.Range("I45").Value = 1556.3334
$ws.Range("K45").Value = 1556.3334
$ws.Range("M45").Value = -1179.3334
$ws.Range("H74").Value = 1162.1957
$ws.Range("I74").Value = 880.34283
$ws.Range("J74").Value = 2059
$ws.Range("K74").Value = 880.34283
$ws.Range("L74").Value = 2059
$ws.Range("M74").Value = -6.342830000000049
$ws.Range("N74").Value = -3807
$ws.Range("H77").Value = 1162.1957
$ws.Range("I77").Value = 880.34283
$ws.Range("J77").Value = 2059
$ws.Range("K77").Value = 4401.71415
$ws.Range("L77").Value = 10295
$ws.Range("M77").Value = -33.71414999999979
$ws.Range("N77").Value = -19031
$ws.Range("H88").Value = 3308.9167
$ws.Range("I88").Value = 1800
$ws.Range("J88").Value = 3446.0908
$ws.Range("K88").Value = 1800
$ws.Range("M88").Value = -1394
$ws.Range("N88").Value = -4258.0908
$ws.Range("H91").Value = 3308.9167
$ws.Range("I91").Value = 1800
$ws.Range("J91").Value = 3446.0908
$ws.Range("K91").Value = 1800
$ws.Range("M91").Value = -396
$ws.Range("N91").Value = -6254.0908
$ws.Range("H110").Value = 1303.9375
$ws.Range("I110").Value = 1030.1482
$ws.Range("K110").Value = 1030.1482
$ws.Range("M110").Value = 1014.8518
$ws.Range("H116").Value = 1389749.8
$ws.Range("I116").Value = 1852333
$ws.Range("K116").Value = 1852333
$ws.Range("M116").Value = -1850039
$ws.Range("H122").Value = 1434.3
$ws.Range("I122").Value = 1397.2174
$ws.Range("J122").Value = 1556.1428
$ws.Range("K122").Value = 4191.6522
$ws.Range("L122").Value = 4668.428400000001
$ws.Range("M122").Value = -1741.6522
$ws.Range("N122").Value = -9568.428400000001

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1389749.8
$ws.Range("I3").Value = 1852333
$ws.Range("K3").Value = 1852333
$ws.Range("M3").Value = -1852219
$ws.Range("H20").Value = 2229.5264
$ws.Range("I20").Value = 2039.75
$ws.Range("K20").Value = 2039.75
$ws.Range("M20").Value = -1792.75
$ws.Range("H99").Value = 1748.5
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1517.5526
$ws.Range("I31").Value = 1109.4445
$ws.Range("K31").Value = 1109.4445
$ws.Range("M31").Value = -814.4445000000001
$ws.Range("H34").Value = 1517.5526
$ws.Range("I34").Value = 1109.4445
$ws.Range("K34").Value = 1109.4445
$ws.Range("M34").Value = -907.4445000000001
$ws.Range("H58").Value = 1892572.8
$ws.Range("J58").Value = 4686.875
$ws.Range("L58").Value = 4686.875
$ws.Range("N58").Value = -5092.875
$ws.Range("H132").Value = 1667.0312
$ws.Range("I132").Value = 1068.909
$ws.Range("J132").Value = 2982.9
$ws.Range("K132").Value = 3206.727
$ws.Range("L132").Value = 8948.700000000001
$ws.Range("M132").Value = -676.7270000000003
$ws.Range("N132").Value = -14008.7
$ws.Range("H134").Value = 1829.6097
$ws.Range("I134").Value = 1810.5927
$ws.Range("K134").Value = 5431.7781
$ws.Range("M134").Value = -2896.7781
$ws.Range("H136").Value = 1892572.8
$ws.Range("J136").Value = 4686.875
$ws.Range("L136").Value = 14060.625
$ws.Range("N136").Value = -19160.625

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 6946.4707
$ws.Range("I56").Value = 6946.4707
$ws.Range("K56").Value = 6946.4707
$ws.Range("M56").Value = -6416.4707

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 45.434784
$ws.Range("I2").Value = 10.076923
$ws.Range("J2").Value = 91.40000000000001
$ws.Range("K2").Value = 10.076923
$ws.Range("L2").Value = 91.40000000000001
$ws.Range("M2").Value = 102.923077
$ws.Range("N2").Value = -317.4
$ws.Range("H19").Value = 55252.75
$ws.Range("I19").Value = 30499.5
$ws.Range("J19").Value = 80006
$ws.Range("K19").Value = 30499.5
$ws.Range("L19").Value = 80006
$ws.Range("M19").Value = -30211.5
$ws.Range("N19").Value = -80582
$ws.Range("H57").Value = 29992.334
$ws.Range("J57").Value = 29992.334
$ws.Range("L57").Value = 29992.334
$ws.Range("N57").Value = -31632.334
$ws.Range("H80").Value = 3582.8333
$ws.Range("I80").Value = 1999.5
$ws.Range("J80").Value = 4374.5
$ws.Range("K80").Value = 1999.5
$ws.Range("L80").Value = 4374.5
$ws.Range("M80").Value = -1001.5
$ws.Range("N80").Value = -6370.5
$ws.Range("H83").Value = 3582.8333
$ws.Range("I83").Value = 1999.5
$ws.Range("J83").Value = 4374.5
$ws.Range("K83").Value = 9997.5
$ws.Range("L83").Value = 21872.5
$ws.Range("M83").Value = -5005.5
$ws.Range("N83").Value = -31856.5
$ws.Range("H113").Value = 1332.8182
$ws.Range("J113").Value = 1333.3334
$ws.Range("L113").Value = 1333.3334
$ws.Range("N113").Value = -5673.3334
$ws.Range("H117").Value = 50310
$ws.Range("J117").Value = 50310
$ws.Range("L117").Value = 50310
$ws.Range("N117").Value = -57194
$ws.Range("H132").Value = 988362.25
$ws.Range("I132").Value = 1375145.5
$ws.Range("J132").Value = 3823.0908
$ws.Range("K132").Value = 4125436.5
$ws.Range("L132").Value = 11469.2724
$ws.Range("M132").Value = -4122906.5
$ws.Range("N132").Value = -16529.2724
$ws.Range("H139").Value = 61082.285
$ws.Range("J139").Value = 61082.285
$ws.Range("L139").Value = 61082.285
$ws.Range("N139").Value = -71362.285

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1011.8823
$ws.Range("I93").Value = 691.8
$ws.Range("J93").Value = 1469.1428
$ws.Range("K93").Value = 691.8
$ws.Range("L93").Value = 1469.1428
$ws.Range("M93").Value = 556.2
$ws.Range("N93").Value = -3965.1428
$ws.Range("H132").Value = 1204.965
$ws.Range("I132").Value = 957
$ws.Range("K132").Value = 2871
$ws.Range("M132").Value = -341
$ws.Range("H136").Value = 2142.8538
$ws.Range("I136").Value = 1370.1
$ws.Range("K136").Value = 4110.299999999999
$ws.Range("M136").Value = -1560.299999999999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H25").Value = 15000
$ws.Range("J25").Value = 15000
$ws.Range("L25").Value = 15000
$ws.Range("N25").Value = -15586
$ws.Range("H132").Value = 1280.8909
$ws.Range("I132").Value = 990.93616
$ws.Range("K132").Value = 2972.80848
$ws.Range("M132").Value = -442.8084799999997
